# Scheduled-runner price/profit refresh for the Leve profit sheets.
# Updates currentAveragePrice(NQ/HQ)/LevePrice(NQ/HQ)/LeveProfit(NQ/HQ)
# columns (H:N) for the rows whose market data changed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 211.14285
$ws.Range("J2").Value = 169.25
$ws.Range("L2").Value = 169.25
$ws.Range("N2").Value = -395.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 770
$ws.Range("J6").Value = 2699
$ws.Range("L6").Value = 8097
$ws.Range("N6").Value = -8321

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 917.1277
$ws.Range("I129").Value = 1286.6666
$ws.Range("K129").Value = 3859.9998
$ws.Range("M129").Value = 1140.0002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 87799
$ws.Range("J133").Value = 87799
$ws.Range("L133").Value = 87799
$ws.Range("N133").Value = -97919

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2880.0754
$ws.Range("I138").Value = 2937.76
$ws.Range("K138").Value = 8813.280000000001
$ws.Range("M138").Value = -3673.280000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 15000479
$ws.Range("I8").Value = 30000000
$ws.Range("J8").Value = 958
$ws.Range("K8").Value = 30000000
$ws.Range("L8").Value = 958
$ws.Range("M8").Value = -29999856
$ws.Range("N8").Value = -1246

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4359.3335
$ws.Range("I32").Value = 3080.7026
$ws.Range("K32").Value = 3080.7026
$ws.Range("M32").Value = -2793.7026

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1506
$ws.Range("I122").Value = 1506
$ws.Range("K122").Value = 4518
$ws.Range("M122").Value = -2068

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1626.8518
$ws.Range("I132").Value = 977.41174
$ws.Range("K132").Value = 2932.23522
$ws.Range("M132").Value = -402.23522

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 120055.18
$ws.Range("I86").Value = 2344
$ws.Range("J86").Value = 288214
$ws.Range("K86").Value = 2344
$ws.Range("L86").Value = 288214
$ws.Range("M86").Value = -1221
$ws.Range("N86").Value = -290460

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 120055.18
$ws.Range("I89").Value = 2344
$ws.Range("J89").Value = 288214
$ws.Range("K89").Value = 11720
$ws.Range("L89").Value = 1441070
$ws.Range("M89").Value = -6104
$ws.Range("N89").Value = -1452302

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2299.3125
$ws.Range("I105").Value = 2306.8462
$ws.Range("J105").Value = 2266.6667
$ws.Range("K105").Value = 2306.8462
$ws.Range("L105").Value = 2266.6667
$ws.Range("M105").Value = -559.8462
$ws.Range("N105").Value = -5760.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 734
$ws.Range("I16").Value = 719.55554
$ws.Range("K16").Value = 719.55554
$ws.Range("M16").Value = -432.55554

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2754.4285
$ws.Range("J31").Value = 4006.2354
$ws.Range("L31").Value = 4006.2354
$ws.Range("N31").Value = -4596.2354

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2754.4285
$ws.Range("J34").Value = 4006.2354
$ws.Range("L34").Value = 4006.2354
$ws.Range("N34").Value = -4410.2354

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2559015.8
$ws.Range("I58").Value = 3624321.2
$ws.Range("K58").Value = 3624321.2
$ws.Range("M58").Value = -3624118.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 734
$ws.Range("I113").Value = 719.55554
$ws.Range("K113").Value = 719.55554
$ws.Range("M113").Value = 1450.44446

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2559015.8
$ws.Range("I136").Value = 3624321.2
$ws.Range("K136").Value = 10872963.6
$ws.Range("M136").Value = -10870413.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 103
$ws.Range("J33").Value = 111
$ws.Range("L33").Value = 666
$ws.Range("N33").Value = -1232

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 275
$ws.Range("I92").Value = 200
$ws.Range("K92").Value = 600
$ws.Range("M92").Value = 648

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 699.8333
$ws.Range("I98").Value = 200
$ws.Range("J98").Value = 799.8
$ws.Range("K98").Value = 600
$ws.Range("L98").Value = 2399.4
$ws.Range("M98").Value = 898
$ws.Range("N98").Value = -5395.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 845.6875
$ws.Range("I107").Value = 250
$ws.Range("J107").Value = 885.4
$ws.Range("K107").Value = 750
$ws.Range("L107").Value = 2656.2
$ws.Range("M107").Value = 1170
$ws.Range("N107").Value = -6496.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 47619470
$ws.Range("I117").Value = 639.5
$ws.Range("K117").Value = 1918.5
$ws.Range("M117").Value = 1523.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1268.5714
$ws.Range("J122").Value = 1230
$ws.Range("L122").Value = 11070
$ws.Range("N122").Value = -15970

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 9960.598
$ws.Range("I131").Value = 575.1111
$ws.Range("J131").Value = 11043.538
$ws.Range("K131").Value = 1725.3333
$ws.Range("L131").Value = 33130.614
$ws.Range("M131").Value = 3314.6667
$ws.Range("N131").Value = -43210.614

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 112.0625
$ws.Range("I2").Value = 230
$ws.Range("J2").Value = 58.454544
$ws.Range("K2").Value = 230
$ws.Range("L2").Value = 58.454544
$ws.Range("M2").Value = -117
$ws.Range("N2").Value = -284.454544

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 2900000
$ws.Range("I14").Value = 2900000
$ws.Range("K14").Value = 2900000
$ws.Range("M14").Value = -2899832

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 19925
$ws.Range("J46").Value = 19925
$ws.Range("L46").Value = 19925
$ws.Range("N46").Value = -20237

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 29999
$ws.Range("J93").Value = 29999
$ws.Range("L93").Value = 29999
$ws.Range("N93").Value = -33743

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3922.7144
$ws.Range("I102").Value = 3922.7144
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3922.7144
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -2300.7144
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1400
$ws.Range("J113").Value = 1500
$ws.Range("L113").Value = 1500
$ws.Range("N113").Value = -5840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1431.75
$ws.Range("I122").Value = 1186.3636
$ws.Range("J122").Value = 1971.6
$ws.Range("K122").Value = 3559.0908
$ws.Range("L122").Value = 5914.799999999999
$ws.Range("M122").Value = -1109.0908
$ws.Range("N122").Value = -10814.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 51073.75
$ws.Range("J139").Value = 51073.75
$ws.Range("L139").Value = 51073.75
$ws.Range("N139").Value = -61353.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3890.3
$ws.Range("I7").Value = 2863
$ws.Range("K7").Value = 2863
$ws.Range("M7").Value = -2751

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3406
$ws.Range("I61").Value = 3140.3333
$ws.Range("K61").Value = 3140.3333
$ws.Range("M61").Value = -2938.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3406
$ws.Range("I113").Value = 3140.3333
$ws.Range("K113").Value = 3140.3333
$ws.Range("M113").Value = -970.3332999999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3890.3
$ws.Range("I126").Value = 2863
$ws.Range("K126").Value = 8589
$ws.Range("M126").Value = -6119

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 15000250
$ws.Range("J6").Value = 500
$ws.Range("L6").Value = 500
$ws.Range("N6").Value = -730

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 580
$ws.Range("I113").Value = 237.5
$ws.Range("K113").Value = 712.5
$ws.Range("M113").Value = 1457.5
